# Apply "Add data for 2022-12-03" update to carjacking-by-month-yoy-latest.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name and workbook sheet reference date
$ws.Name = "Through 2022-11-25"

# Update the November row label (shared string used in A12)
$ws.Range("A12").Value = "November (through 11-25)"

# Update November (row 12) figures for columns C..I (2016..2022)
$ws.Range("C12").Value = 62
$ws.Range("D12").Value = 94
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 45
$ws.Range("G12").Value = 180
$ws.Range("H12").Value = 170
$ws.Range("I12").Value = 92

# Update Total row (row 13) figures for columns C..I (2016..2022)
$ws.Range("C13").Value = 548
$ws.Range("D13").Value = 804
$ws.Range("E13").Value = 665
$ws.Range("F13").Value = 527
$ws.Range("G13").Value = 1237
$ws.Range("H13").Value = 1611
$ws.Range("I13").Value = 1490
